$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NATMI ligand-receptor pair metrics (rows 2-7, columns E:T)
# following the re-run of the analysis with the revised parameters.
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.270036666666666
$ws.Cells.Item(2, 8).Value = 9.81011
$ws.Cells.Item(2, 9).Value = 0.359406393324744
$ws.Cells.Item(2, 10).Value = 0.3594063933247441
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 62.90731233333333
$ws.Cells.Item(2, 14).Value = 188.721937
$ws.Cells.Item(2, 15).Value = 0.9393635410440488
$ws.Cells.Item(2, 16).Value = 0.9393635410440487
$ws.Cells.Item(2, 17).Value = 205.7092179314522
$ws.Cells.Item(2, 18).Value = 1851.38296138307
$ws.Cells.Item(2, 19).Value = 0.3376132623074017
$ws.Cells.Item(2, 20).Value = 0.3376132623074017
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.270036666666666
$ws.Cells.Item(3, 8).Value = 9.81011
$ws.Cells.Item(3, 9).Value = 0.359406393324744
$ws.Cells.Item(3, 10).Value = 0.3594063933247441
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.253965666666666
$ws.Cells.Item(3, 14).Value = 6.761896999999999
$ws.Cells.Item(3, 15).Value = 0.03365734588711396
$ws.Cells.Item(3, 16).Value = 0.03365734588711396
$ws.Cells.Item(3, 17).Value = 7.370550375407777
$ws.Cells.Item(3, 18).Value = 66.33495337867
$ws.Cells.Item(3, 19).Value = 0.01209666529417103
$ws.Cells.Item(3, 20).Value = 0.01209666529417104
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.270036666666666
$ws.Cells.Item(4, 8).Value = 9.81011
$ws.Cells.Item(4, 9).Value = 0.359406393324744
$ws.Cells.Item(4, 10).Value = 0.3594063933247441
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.806737666666667
$ws.Cells.Item(4, 14).Value = 5.420213
$ws.Cells.Item(4, 15).Value = 0.02697911306883729
$ws.Cells.Item(4, 16).Value = 0.02697911306883729
$ws.Cells.Item(4, 17).Value = 5.908098417047777
$ws.Cells.Item(4, 18).Value = 53.17288575343
$ws.Cells.Item(4, 19).Value = 0.009696465723171276
$ws.Cells.Item(4, 20).Value = 0.009696465723171278
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.828401
$ws.Cells.Item(5, 8).Value = 17.485203
$ws.Cells.Item(5, 9).Value = 0.6405936066752559
$ws.Cells.Item(5, 10).Value = 0.640593606675256
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 62.90731233333333
$ws.Cells.Item(5, 14).Value = 188.721937
$ws.Cells.Item(5, 15).Value = 0.9393635410440488
$ws.Cells.Item(5, 16).Value = 0.9393635410440487
$ws.Cells.Item(5, 17).Value = 366.6490421109123
$ws.Cells.Item(5, 18).Value = 3299.841378998211
$ws.Cells.Item(5, 19).Value = 0.601750278736647
$ws.Cells.Item(5, 20).Value = 0.6017502787366471
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5.828401
$ws.Cells.Item(6, 8).Value = 17.485203
$ws.Cells.Item(6, 9).Value = 0.6405936066752559
$ws.Cells.Item(6, 10).Value = 0.640593606675256
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.253965666666666
$ws.Cells.Item(6, 14).Value = 6.761896999999999
$ws.Cells.Item(6, 15).Value = 0.03365734588711396
$ws.Cells.Item(6, 16).Value = 0.03365734588711396
$ws.Cells.Item(6, 17).Value = 13.13701574556567
$ws.Cells.Item(6, 18).Value = 118.233141710091
$ws.Cells.Item(6, 19).Value = 0.02156068059294292
$ws.Cells.Item(6, 20).Value = 0.02156068059294292
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5.828401
$ws.Cells.Item(7, 8).Value = 17.485203
$ws.Cells.Item(7, 9).Value = 0.6405936066752559
$ws.Cells.Item(7, 10).Value = 0.640593606675256
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.806737666666667
$ws.Cells.Item(7, 14).Value = 5.420213
$ws.Cells.Item(7, 15).Value = 0.02697911306883729
$ws.Cells.Item(7, 16).Value = 0.02697911306883729
$ws.Cells.Item(7, 17).Value = 10.53039162313767
$ws.Cells.Item(7, 18).Value = 94.77352460823901
$ws.Cells.Item(7, 19).Value = 0.01728264734566601
$ws.Cells.Item(7, 20).Value = 0.01728264734566601
